$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell content updates (Testeo/Plan de Testo.xlsx, sheet "Wildo") ---
$ws.Range('D91').Value = 'Verificado, aparece el numero de la nota de retiro en el autocomplete.'
$ws.Range('D93').Value = 'Aparecen los destinatarios al ingresar "2" letras en el autocomplete.'
$ws.Range('E93').Value = 161
$ws.Range('D94').Value = 'Verificado, al seleccionar el focus va sobre código de barras'
$ws.Range('D96').Value = 'No verificado'
$ws.Range('D97').Value = 'Verificado, no se puede guardar mas que "Cantidad Productos" productos'
$ws.Range('D99').Value = 'Verificado, aparece los productos registrados anteriormente.'
$ws.Range('D101').Value = 'Verificado, es posible registrar un nuevo destinatario y luego usarlo'
$ws.Range('D102').Value = 'Verificado.'
$ws.Range('D103').Value = 'Verificado, se guardan correctamente.'
$ws.Range('D105').Value = 'Verificado.'
$ws.Range('D106').Value = 'Verificado.'
$ws.Range('D107').Value = 'Verificado, no permite ingresar 2 codigos de barra iguales'
$ws.Range('D109').Value = 'Verificado.'
$ws.Range('D110').Value = 'Verificado.'
$ws.Range('D111').Value = 'Verificado.'
$ws.Range('D112').Value = 'Verificado, se comporta como esperado.'
$ws.Range('B113').Value = 'Otros'
$ws.Range('C113').Value = 'Modificar fecha'
$ws.Range('D113').Value = 'Al modificar la fecha no guarda el producto'
$ws.Range('E113').Value = 160
$ws.Range('B114').Value = 'Deberia guardar el producto al tener una direccion valida'
$ws.Range('C114').Value = 'Al seleccionar la dirección "seleccionar", hacer click'
$ws.Range('D114').Value = 'No guarda el producto.'
$ws.Range('E114').Value = 162
$ws.Range('C115').Value = 'en guardar, luego intentar seleccionar una dirección valida'

# C116 changes text AND loses its bold styling (reverts to Normal style)
$ws.Range('C116').Value = 'y hacer click en guardar'
$ws.Range('C116').Style = 'Normal'

# --- Sheet view: selection moved to D115; scroll position updated to C90 ---
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 90
    $win.ScrollColumn = 3
} catch {}
$ws.Range('D115').Select()
